# 农用薄膜.xlsx edit
#
# Within every 4-row "year block" (the sub-periods A/B/C/D stacked one
# after another starting at row 2), the "B" sub-period row and the "C"
# sub-period row swap their entire A:E content - i.e. the row that used
# to carry the "...年B" label/data now carries the "...年C" label/data,
# and vice versa. (Rows A and D in each block are untouched.)
#
# Afterwards, columns F ("农用薄膜产销率") and G ("农用薄膜销售量"),
# including their header cells in row 1, are deleted outright, which also
# shrinks the sheet's used range from A1:G65 down to A1:E65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (B-row, C-row) pairs for each of the 16 year blocks (2000 .. 2015).
$pairs = @(
    @(3,4),   @(7,8),   @(11,12), @(15,16),
    @(19,20), @(23,24), @(27,28), @(31,32),
    @(35,36), @(39,40), @(43,44), @(47,48),
    @(51,52), @(55,56), @(59,60), @(63,64)
)

foreach ($pair in $pairs) {
    $rB = $pair[0]
    $rC = $pair[1]

    foreach ($col in @('A','B','C','D','E')) {
        $cellB = $ws.Range("$col$rB")
        $cellC = $ws.Range("$col$rC")

        $valB = $cellB.Value2
        $valC = $cellC.Value2

        if ($valB -ne $valC) {
            $cellB.Value = $valC
            $cellC.Value = $valB
        }
    }
}

# Drop the F and G columns (values + headers) entirely.
$ws.Range("F1:G65").Delete()
